$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Comment.Delete()
$ws.Columns("B").Delete()
$ws.Range("I2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:circlek@yahoo.com")
$ws.Range("H2").Style = "Hyperlink"
$ws.Range("K4").Select() | Out-Null
